# Modify full search method
# Update rule-search result table (rows 1-20) with new metric values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rule1
$ws.Range("B1").Value = 110
$ws.Range("C1").Value = 110
$ws.Range("D1").Value = 400
$ws.Range("E1").Value = 550
$ws.Range("F1").Value = 70
$ws.Range("G1").Value = 200
$ws.Range("H1").Value = 10
$ws.Range("I1").Value = 0.2633997766844854

# rule2
$ws.Range("B2").Value = 110
$ws.Range("C2").Value = 110
$ws.Range("D2").Value = 400
$ws.Range("E2").Value = 550
$ws.Range("F2").Value = 70
$ws.Range("G2").Value = 200
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 0.2633997766844854

# rule3
$ws.Range("B3").Value = 110
$ws.Range("C3").Value = 110
$ws.Range("D3").Value = 400
$ws.Range("E3").Value = 550
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = 200
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 0.2633997766844854

# rule4
$ws.Range("B4").Value = 110
$ws.Range("C4").Value = 110
$ws.Range("D4").Value = 400
$ws.Range("E4").Value = 550
$ws.Range("F4").Value = 70
$ws.Range("G4").Value = 160
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 0.2633997766844854

# rule5
$ws.Range("B5").Value = 110
$ws.Range("C5").Value = 110
$ws.Range("D5").Value = 400
$ws.Range("E5").Value = 550
$ws.Range("F5").Value = 70
$ws.Range("G5").Value = 200
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 0.2633997766844854

# rule6
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 110
$ws.Range("D6").Value = 300
$ws.Range("E6").Value = 550
$ws.Range("F6").Value = 70
$ws.Range("G6").Value = 200
$ws.Range("H6").Value = "NOT ADAPTED"
$ws.Range("I6").Value = -0.1364463133780686

# rule7
$ws.Range("B7").Value = 100
$ws.Range("C7").Value = 110
$ws.Range("D7").Value = 400
$ws.Range("E7").Value = 550
$ws.Range("F7").Value = 70
$ws.Range("G7").Value = 200
$ws.Range("H7").Value = "NOT ADAPTED"
$ws.Range("I7").Value = -0.1364463133780686

# rule8
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 110
$ws.Range("D8").Value = 400
$ws.Range("E8").Value = 1100
$ws.Range("F8").Value = 70
$ws.Range("G8").Value = 200
$ws.Range("H8").Value = "NOT ADAPTED"
$ws.Range("I8").Value = -0.1364463133780686

# rule9
$ws.Range("B9").Value = 60
$ws.Range("C9").Value = 110
$ws.Range("D9").Value = 550
$ws.Range("E9").Value = 1100
$ws.Range("F9").Value = 70
$ws.Range("G9").Value = 200
$ws.Range("H9").Value = "NOT ADAPTED"
$ws.Range("I9").Value = -0.1364463133780686

# rule10
$ws.Range("B10").Value = 100
$ws.Range("C10").Value = 110
$ws.Range("D10").Value = 400
$ws.Range("E10").Value = 550
$ws.Range("F10").Value = 70
$ws.Range("G10").Value = 200
$ws.Range("H10").Value = "NOT ADAPTED"
$ws.Range("I10").Value = -0.1364463133780686

# rule11
$ws.Range("B11").Value = 100
$ws.Range("C11").Value = 110
$ws.Range("D11").Value = 400
$ws.Range("E11").Value = 1100
$ws.Range("F11").Value = 70
$ws.Range("G11").Value = 200
$ws.Range("H11").Value = "NOT ADAPTED"
$ws.Range("I11").Value = -0.1364463133780686

# rule12
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = 110
$ws.Range("D12").Value = 400
$ws.Range("E12").Value = 1100
$ws.Range("F12").Value = 70
$ws.Range("G12").Value = 200
$ws.Range("H12").Value = "NOT ADAPTED"
$ws.Range("I12").Value = -0.1364463133780686

# rule13
$ws.Range("B13").Value = 100
$ws.Range("C13").Value = 110
$ws.Range("D13").Value = 400
$ws.Range("E13").Value = 550
$ws.Range("F13").Value = 70
$ws.Range("G13").Value = 200
$ws.Range("H13").Value = "NOT ADAPTED"
$ws.Range("I13").Value = -0.1364463133780686

# rule14
$ws.Range("B14").Value = 100
$ws.Range("C14").Value = 110
$ws.Range("D14").Value = 400
$ws.Range("E14").Value = 550
$ws.Range("F14").Value = 70
$ws.Range("G14").Value = 200
$ws.Range("H14").Value = "NOT ADAPTED"
$ws.Range("I14").Value = -0.1364463133780686

# rule15
$ws.Range("B15").Value = 100
$ws.Range("C15").Value = 110
$ws.Range("D15").Value = 400
$ws.Range("E15").Value = 550
$ws.Range("F15").Value = 70
$ws.Range("G15").Value = 200
$ws.Range("H15").Value = "NOT ADAPTED"
$ws.Range("I15").Value = -0.1364463133780686

# rule16
$ws.Range("B16").Value = 100
$ws.Range("C16").Value = 110
$ws.Range("D16").Value = 400
$ws.Range("E16").Value = 1100
$ws.Range("F16").Value = 70
$ws.Range("G16").Value = 200
$ws.Range("H16").Value = "NOT ADAPTED"
$ws.Range("I16").Value = -0.1364463133780686

# rule17
$ws.Range("B17").Value = 100
$ws.Range("C17").Value = 110
$ws.Range("D17").Value = 300
$ws.Range("E17").Value = 550
$ws.Range("F17").Value = 70
$ws.Range("G17").Value = 200
$ws.Range("H17").Value = "NOT ADAPTED"
$ws.Range("I17").Value = -0.1364463133780686

# rule18
$ws.Range("B18").Value = 100
$ws.Range("C18").Value = 110
$ws.Range("D18").Value = 400
$ws.Range("E18").Value = 550
$ws.Range("F18").Value = 70
$ws.Range("G18").Value = 200
$ws.Range("H18").Value = "NOT ADAPTED"
$ws.Range("I18").Value = -0.1364463133780686

# rule19
$ws.Range("B19").Value = 100
$ws.Range("C19").Value = 110
$ws.Range("D19").Value = 400
$ws.Range("E19").Value = 550
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 200
$ws.Range("H19").Value = "NOT ADAPTED"
$ws.Range("I19").Value = -0.1364463133780686

# rule20
$ws.Range("B20").Value = 100
$ws.Range("C20").Value = 110
$ws.Range("D20").Value = 400
$ws.Range("E20").Value = 800
$ws.Range("F20").Value = 70
$ws.Range("G20").Value = 140
$ws.Range("H20").Value = "NOT ADAPTED"
$ws.Range("I20").Value = -0.1364463133780686

